$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.767.19'
$ws.Range("E2").Value = '  +0.81%  '

$ws.Range("D3").Value = '1.648.11'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("E4").Value = '  +0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '

$ws.Range("E7").Value = '  +0.41%  '

$ws.Range("E8").Value = '  +0.72%  '

$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("D12").Value = '1.876.42'
$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("D13").Value = '1.657.92'
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.29%  '

$ws.Range("E15").Value = '  +1.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '26.748.97'
$ws.Range("E17").Value = '  +0.67%  '

$ws.Range("D18").Value = '0.0₃0744'
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.28%  '

$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +16.10%  '

$ws.Range("E22").Value = '  +1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("E28").Value = '  +4.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.21%  '

$ws.Range("E30").Value = '  +1.24%  '

$ws.Range("E31").Value = '  +1.61%  '

$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("E33").Value = '  +1.23%  '

$ws.Range("D34").Value = '1.279.72'
$ws.Range("E34").Value = '  +0.98%  '

$ws.Range("E35").Value = '  +2.59%  '

$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.541'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.24%  '

$ws.Range("E40").Value = '  +0.44%  '

$ws.Range("E41").Value = '  +2.40%  '

$ws.Range("E42").Value = '  -0.79%  '

$ws.Range("E43").Value = '  +1.71%  '

$ws.Range("D44").Value = '1.787.96'
$ws.Range("E44").Value = '  +0.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.05%  '

$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  +1.46%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0515'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.74%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0979'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.30%  '
